# Update "想去人数" (F column) values for several events across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 632
$wsExpo.Range("F5").Value = 555
$wsExpo.Range("F6").Value = 306
$wsExpo.Range("F7").Value = 2766
$wsExpo.Range("F8").Value = 469
$wsExpo.Range("F9").Value = 7738
$wsExpo.Range("F10").Value = 200
$wsExpo.Range("F12").Value = 40
$wsExpo.Range("F13").Value = 322
$wsExpo.Range("F14").Value = 47

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 4

# --- Sheet "全部类型" (All types, combined static copy) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 632
$wsAll.Range("F5").Value = 555
$wsAll.Range("F6").Value = 306
$wsAll.Range("F9").Value = 2766
$wsAll.Range("F10").Value = 469
$wsAll.Range("F11").Value = 7738
$wsAll.Range("F12").Value = 200
$wsAll.Range("F14").Value = 40
$wsAll.Range("F15").Value = 4
$wsAll.Range("F17").Value = 322
$wsAll.Range("F18").Value = 47
